$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 8
$ws.Range("B113").Value = "Terminal La Palmera de La Serena"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44680
$ws.Range("E113").Value = 4
$ws.Range("F113").Value = 100112037
$ws.Range("G113").Value = "Cebollín"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1150
$ws.Range("N113").Value = '$/paquete 6 unidades'
$ws.Range("O113").Value = "Provincia del Elquí"
$ws.Range("P113").Value = 192
$ws.Range("Q113").Value = 6
$ws.Range("R113").Value = "Hortaliza"
